$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the SqlServer-section (TestFile column F, rows 20/21/24/25) labels
# to match the InMemory-section (rows 5/6/9/10) naming so that the
# CreateInput/UpdateInput rows use the plural "Inputs" form consistently.
$ws.Range("F20").Value = "CreateInputs0"
$ws.Range("F21").Value = "CreateInputs1"
$ws.Range("F24").Value = "UpdateInputs0"
$ws.Range("F25").Value = "UpdateInputs1"

# Update the view state: scroll position and selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("F17:F31").Select()

$wb.Save()
